# Add data for 2022-08-29 (Chicago carjacking-by-neighborhood-by-month workbook)
# Updates the "through August 20" snapshot to "through August 21" and adds
# the incremental carjacking counts recorded on 2022-08-29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "Through 2022-08-21"

# Update the header label for column B (also updates the shared string)
$ws.Range("B1").Value = "August 2022 (through August 21)"

# Update existing cell values (row => column => new value)
$ws.Range("AP3").Value = 3
$ws.Range("AX3").Value = 4

$ws.Range("J5").Value = 10
$ws.Range("AH5").Value = 3
$ws.Range("AP5").Value = 2

$ws.Range("B6").Value = 7

$ws.Range("R7").Value = 5

$ws.Range("AX8").Value = 2

$ws.Range("R12").Value = 5
$ws.Range("AX12").Value = 2

$ws.Range("B18").Value = 2

$ws.Range("R19").Value = 3

$ws.Range("J21").Value = 2

$ws.Range("J24").Value = 2

# New cells that previously had no value
$ws.Range("B28").Value = 1
$ws.Range("B45").Value = 1

$ws.Range("J45").Value = 3

$ws.Range("R46").Value = 2
